$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws2 = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" headers to the new metric-specific names
$ws1.Range("B1").Value = "Weekly_PO_Qty"
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet, placed after "Monthly Trend"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# Header row
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# Copy the header formatting (bold, centered, bordered) from the existing sheets
$ws1.Range("B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Forecast data rows
$data = @(
    @(45242.99999999999, 72,  -31.76756608461403,  170.0568341045128),
    @(45277.99999999999, 97,   -4.240868695654157, 190.7474990394152),
    @(45298.99999999999, 112,   7.217988793849728, 210.2977667299166),
    @(45312.99999999999, 122,  22.07366371618866,  216.3743018688133),
    @(45319.99999999999, 127,  25.19186171036166,  227.0928036411529),
    @(45326.99999999999, 132,  31.54466270261188,  232.5940144352242),
    @(45333.99999999999, 137,  35.6833564892454,   239.2203612118289),
    @(45340.99999999999, 142,  43.68986885525967,  240.0515699028344),
    @(45347.99999999999, 147,  39.12479491894153,  246.8446709757068),
    @(45354.99999999999, 153,  54.20289978780182,  260.4081525675414),
    @(45361.99999999999, 158,  63.51831085392725,  258.4203716305223),
    @(45368.99999999999, 163,  65.26969408152669,  264.3614771643931),
    @(45375.99999999999, 168,  67.96042495121195,  261.9981817611309),
    @(45382.99999999999, 173,  64.8438758861471,   271.6408429529968),
    @(45389.99999999999, 178,  78.57830168534699,  276.4950467743939)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws3.Cells.Item($row, 1).Value = $data[$i][0]
    $ws3.Cells.Item($row, 2).Value = $data[$i][1]
    $ws3.Cells.Item($row, 3).Value = $data[$i][2]
    $ws3.Cells.Item($row, 4).Value = $data[$i][3]
}

# Apply the date-formatted style used in column A of the other sheets
$ws1.Range("A2").Copy()
$ws3.Range("A2:A16").PasteSpecial(-4122)
